$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values first ---

# Row 8: leo palette
$ws.Range("A8").Value = "leo"
$ws.Range("B8").Value = "#0081A7"
$ws.Range("C8").Value = "#000000"
$ws.Range("D8").Value = "#FEFEF1"
$ws.Range("E8").Value = "#000000"
$ws.Range("F8").Value = "#004154, #00556F, #006D8D, #0087A9, #00A8B6, #48C5C3, #A2E0CF, #FDFCDC, #FDE3C1, #F8AC94, #EA6E64, #CC5F58, #994C47, #6D3C39, #523332"
$ws.Range("G8").Value = "#0081A7, #f07167, #A2E0CF, #FDE3C1"

# Row 9: portal palette
$ws.Range("A9").Value = "portal"
$ws.Range("B9").Value = "#0B2B51"
$ws.Range("C9").Value = "#061629"
$ws.Range("D9").Value = "#FEFEF1"
$ws.Range("E9").Value = "#F3E4C2"
$ws.Range("F9").Value = "#0B2B51, #0E3768, #124482, #16549E, #5B8EBA, #A5C7D0, #ECF2DD, #F9F0CD, #F2DAAF, #DEAA79, #C88F60, #AD7C54, #815F43, #614A38, #4C3D32"
$ws.Range("G9").Value = "#0E3564, #D89B68, #815F43, #A5C7D0"

# Row 10: pem palette
$ws.Range("A10").Value = "pem"
$ws.Range("B10").Value = "#4C0C12"
$ws.Range("C10").Value = "#1F1F1F"
$ws.Range("D10").Value = "#FFFBEB"
$ws.Range("E10").Value = "#4C0C12"
$ws.Range("F10").Value = "#112438, #342031, #571C2A, #7A1823, #992E2E, #B5604B, #D19168, #EDC385, #DBA563, #C98841, #B76B1F, #9C5710, #764C14, #504118, #2A361C"
$ws.Range("G10").Value = "#0D1B2A, #AF5D0F, #DEB77D, #354322"

# --- Formatting ---
# G9 gets a brand-new numeric-format style.
$ws.Range("G9").NumberFormat = "0.00E+00"

# G8 / G10 reuse the existing "Helvetica" cell style already used by column F
# (e.g. F2) - copy its format instead of touching Font.Name directly, which
# would otherwise register a redundant duplicate font in the style table.
$ws.Range("F2").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("G10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the current selection to match the saved worksheet view.
$ws.Range("F17").Select() | Out-Null
